$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace row 9 (previously the "ABM7" 16MHz crystal) with the new "NX3225"
# 16MHz crystal part, reusing the description/manufacturer/part-number/
# digikey strings already used by the NDK NX3225SA crystal row (row 23).
$ws.Range("B9").Value = "NX3225"
$ws.Range("C9").Value = "4-SMD"
$ws.Range("F9").Value = "CRYSTAL 16MHZ 7.2PF SMD"
$ws.Range("G9").Value = "NDK"
$ws.Range("H9").Value = "NX3225SA-16.000000MHZ-B3"
$ws.Range("I9").Value = "644-1099-1-ND"

$ws.Range("K9").Value = 0.95
$ws.Range("L9").Value = 0.84
$ws.Range("M9").Value = 0.69
$ws.Range("N9").Value = 0.54

$ws.Range("B10").Select()
